# Update machine_gap_summary: remove setup_efficiency column (G) and
# adjust setup/idle/gap time values reflecting fixed DelayProcessor integration.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update changed numeric values first (before removing the column,
# so row/column references stay simple).
$ws.Range("D2").Value = 126
$ws.Range("E2").Value = 59

$ws.Range("C5").Value = 294
$ws.Range("D5").Value = 218
$ws.Range("E5").Value = 76

# Remove the now-obsolete "setup_efficiency" column entirely.
$ws.Range("G1:G6").Delete()
